$d = $word.ActiveDocument

# The "invitacion y sorteo" placeholder row ([[#productos]] ... [[/productos]])
# in the first table currently has single-underline formatting on every run.
# Remove the underline so the merge-field row renders as plain text.
$table = $d.Tables.Item(1)
$row = $table.Rows.Item(2)
$row.Range.Font.Underline = 0
